# Fix merge conflicts and update scenario names
# - Rename the 4 existing scenario sheets, add 2 new ones (Moderate-Mid (SC), Expanded-High)
# - Fix the "Cummulative Capacity" -> "Cumulative Capacity" header typo on every sheet
# - Replace the Year / Cumulative Capacity data series on every sheet with the corrected numbers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the four pre-existing sheets (order/sheetId stays the same, only
#    the display name changes) and fix the old "Cummulative Capacity" typo.
# ---------------------------------------------------------------------------
$newNames = @("Baseline-Low", "Baseline-Mid (SC)", "Baseline-Mid (CC)", "Moderate-Low")
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
    $ws.Range("B1").Value = "Cumulative Capacity"
}

# ---------------------------------------------------------------------------
# 2. Add the two brand-new scenario sheets at the end of the workbook, copying
#    the header formatting (bold, centered, bordered) from the first sheet so
#    the new tabs match the existing look.
# ---------------------------------------------------------------------------
$headerSrc = $wb.Worksheets.Item(1)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5 = $wb.Worksheets.Add($null, $lastSheet)
$sheet5.Name = "Moderate-Mid (SC)"
$headerSrc.Range("A1:B1").Copy()
$sheet5.Range("A1:B1").PasteSpecial(-4122)
$sheet5.Range("A1").Value = "Year"
$sheet5.Range("B1").Value = "Cumulative Capacity"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet6 = $wb.Worksheets.Add($null, $lastSheet)
$sheet6.Name = "Expanded-High"
$headerSrc.Range("A1:B1").Copy()
$sheet6.Range("A1:B1").PasteSpecial(-4122)
$sheet6.Range("A1").Value = "Year"
$sheet6.Range("B1").Value = "Cumulative Capacity"

# ---------------------------------------------------------------------------
# 3. Data series (Year in column A, Cumulative Capacity in column B) for each
#    of the 6 sheets, in final tab order.
# ---------------------------------------------------------------------------
$years1 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049,2050,2051)
$vals1  = @(1298.617392332037,2276.697484881295,3261.977426753822,4242.837463853284,5218.916186298174,6190.215016890826,7163.164940367873,8214.825658884853,9520.265118171012,10822.57966268448,12153.38094730442,13478.5894208557,14786.29677183206,16095.7060532272,17437.07746674893,18749.47402992389,20061.73197624848,21396.43398280865,22711.40342408068,23999.28333973481,24872.50337888572)

$years2 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046)
$vals2  = @(1154.17065076776,1823.443388267057,2817.795556914285,4127.643693473954,6195.251873359175,8343.191712963724,10474.009708605,12606.2351429245,14747.74403663507,16442.10526315789,17880.86069368357,19442.90448930401,20880.86069368357,22442.90448930401,23546.14203467721,24748.92988929889)

$years3 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049,2050,2051,2052)
$vals3  = @(1154.17065076776,1823.443388267057,2331.174317093777,3155.730780061955,4250.915979584492,5426.359516616314,6442.105263157895,7928.690832229717,9415.276401301538,10901.83643400882,12390.5272398039,13873.61692525089,15356.2881468058,16830.27245653833,18406.88788507601,20204.79053293633,21101.2252042007,21740.08168028005,22393.24996233238,23051.34918903193,23688.48643537712,24338.09092282133)

$years4 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049,2050,2051)
$vals4  = @(1298.617392332037,2609.78534530162,4092.65776437349,5572.232441022234,7048.311163467123,8519.609994059774,9993.929780550521,11665.62301399959,14312.00153677913,16942.74542067324,19634.15229886452,22285.46292307787,24567.28604477605,26095.7060532272,27437.07746674893,28749.47402992389,30061.73197624848,31396.43398280865,32711.40342408068,33999.28333973481,34872.50337888572)

$years5 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048)
$vals5  = @(1154.17065076776,2156.531248687382,3648.475894533953,5457.038670642904,7538.709090501497,9700.711170079419,11702.26015237906,13804.96388886606,16619.91626680589,19420.43010715657,22377.76500531643,25349.25290401158,27962.22515081275,29832.53836002874,31595.74564221231,33435.67074317917,34336.27850777828,34997.80999848965)

$years6 = @(2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049,2050,2051)
$vals6  = @(1154.17065076776,2156.531248687382,3648.475894533953,5457.038670642904,7538.709090501497,9700.711170079419,11702.26015237906,13804.96388886606,16619.91626680589,20304.64063347235,24579.09341860057,28843.57024718502,33106.59742401617,37420.80771006332,41796.22882817144,46262.73980836057,49818.24025444769,52309.18932698269,53330.86419753086,53986.30751964085,54630.01322168356)

$allYears = @($years1, $years2, $years3, $years4, $years5, $years6)
$allVals  = @($vals1,  $vals2,  $vals3,  $vals4,  $vals5,  $vals6)

for ($s = 1; $s -le 6; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $years = $allYears[$s - 1]
    $vals = $allVals[$s - 1]

    for ($i = 0; $i -lt $years.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $years[$i]
        $ws.Cells.Item($row, 2).Value = $vals[$i]
    }

    # Drop any leftover rows below the new data range (sheet1 shrank from 26
    # to 22 rows; sheets 3/4 grew, so this is a no-op range for those).
    $lastRow = $years.Length + 1
    $ws.Range("A" + ($lastRow + 1) + ":B30").ClearContents()
}

# Keep the first tab active/selected, matching the unchanged activeTab="0".
$wb.Worksheets.Item(1).Activate()
